$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo'd domain in C9 (jasurgraduate -> jasurlive)
$ws.Range("C9").Value = "I created a quiz of testing your Hangul knowledge. And in the end according to your result you can get a personal certificate like this. Give it a try! You will love it! | https://jasurlive.github.io/Hangul/"

# Replace stale image link in D11 with a fresh one
$ws.Range("D11").Value = "https://blogger.googleusercontent.com/img/a/AVvXsEhRneIHityZCfdxYN2EabHWzPEHpiqWFw9UFIPEzeBZNwjGyDQK-M4bQ2ZMCA8SdgZ_k1UYS-eKWhrU3uF_V87SKRyqE7Fi-TCW11UoEfGRvH-zfygoCGuV5hCOMyMA9Ty-Xzj9AjW5C7_B255wIj_ZSHP52H9ExRvatwOquMS-B99GUUv7_0x5k9TJ8p8"

# Add a new portfolio entry row 12
$ws.Range("A12").Value = "2025-05-22 17:34:06"
$ws.Range("B12").Value = "POS App for a local market in Daejeon"
$ws.Range("C12").Value = "Point Of Sale (POS) app for the local store. Inventory control, calculations of income and outcome. Works with any inventory. The excel sheet needs to be modified accordingly before using it."
$ws.Range("D12").Value = "https://blogger.googleusercontent.com/img/b/R29vZ2xl/AVvXsEjA1FzCHoSKmitLjioCBqb4eMmq3MQAFdEnUhzNt5q6WthEYhnTavPmURhhQFcGI46EKgr9SokoE00hfF87GYdCmjUk3YZGgBMO6HW4V8t_tpgg2UHF0rZnVl8Df15AXWG7kZbBLCBQvlwbYBDHhZ3tPJDBGAWmwyVcI-UzfstOXxwSAb2UnIaet9xpYUQ/s1917/POS.png"
